$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 12.94719999999999
$ws.Range("E8").Value = 13.33389999999999
$ws.Range("B12").Value = 5.357000000000002
$ws.Range("E12").Value = 12.54189999999999
$ws.Range("E14").Value = 13.7127
$ws.Range("E22").Value = 11.9342
